$d = $word.ActiveDocument

# --- Remove the blank paragraph between ${group} and the table ------
$blankPara = $d.Paragraphs.Item(2)
$blankPara.Range.Delete()

# --- Move the _GoBack bookmark from the signDate paragraph to the ---
# --- ${group} paragraph ---------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$groupPara = $d.Paragraphs.Item(1)
$groupRange = $d.Range($groupPara.Range.Start, $groupPara.Range.End)
$d.Bookmarks.Add("_GoBack", $groupRange)

# --- Table changes -------------------------------------------------
# The table starts as 1 row x 3 columns: [${!fullName}] [${!birthDate}] [empty]
# Target: 2 rows x 2 columns:
#   Row 1: [Full name] [Birth date]            (new header-style row)
#   Row 2: [${!fullName}] [${!birthDate}]      (the original placeholder row, shifted down)
$t = $d.Tables.Item(1)
$firstRow = $t.Rows.Item(1)

# Insert a brand-new empty row above the existing one; the existing row
# (with its original placeholder runs/proofErr markup intact) becomes row 2.
$t.Rows.Add($firstRow)

# Drop the now-unused third column (was only ever an empty spacer cell).
$t.Columns.Item(3).Delete()

# Fill in the new header row with plain labels.
$t.Cell(1, 1).Range.InsertAfter("Full name")
$t.Cell(1, 2).Range.InsertAfter("Birth date")
